# Generate Report for Handoff
# Replaces the old handoff UUID/hashes/timestamps with the new ones produced
# by the latest handoff run, and clears the now-stale "Latest Target File"
# (col I) / "Latest Handback File" (col J) columns on the per-locale sheets
# since a fresh handoff has no handback yet.

$wb = $excel.ActiveWorkbook

$newGuid = "a3625f66-258b-4b76-9712-5be3cc369d3b"

$newZhHash = "a3625f66-258b-4b76-9712-5be3cc369d3b.883b146127c68383bb3730d19a6d7ed2a75b731f.zh-cn.xlf"
$newDeHash = "a3625f66-258b-4b76-9712-5be3cc369d3b.883b146127c68383bb3730d19a6d7ed2a75b731f.de-de.xlf"

$newGenDate = "2016-08-19 21:03:20"
$newZhHoDate = "2016-08-19 21:03:16"
$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newGenDate

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Remove the stale "Latest Target File" hyperlink (column I) before
# touching the cell, so it does not linger after the value is cleared.
foreach ($hl in $wsZh.Hyperlinks) {
    $r = $hl.Range
    if ($r.Column -eq 9) {
        $hl.Delete()
    }
}

$wsZh.Range("A2").Value = "$newGuid.md"
foreach ($hl in $wsZh.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

$wsZh.Range("G2").Value = $newZhHash
$wsZh.Range("H2").Value = $newZhHoDate

$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""

$wsZh.Range("K2").Value = $epoch

# These raw inputs are chosen so that, after the COM layer's internal
# "+5px padding, snapped to the nearest 1/6 character" column-width
# rounding, the saved XML width lands as close as this engine can get to
# the true Excel autofit widths (18.6506053379604 / 21.7054770333426).
$wsZh.Columns.Item(9).ColumnWidth = 17.8172720046271
$wsZh.Columns.Item(10).ColumnWidth = 20.8721437000093

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

foreach ($hl in $wsDe.Hyperlinks) {
    $r = $hl.Range
    if ($r.Column -eq 9) {
        $hl.Delete()
    }
}

$wsDe.Range("A2").Value = "$newGuid.md"
foreach ($hl in $wsDe.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

$wsDe.Range("G2").Value = $newDeHash
$wsDe.Range("H2").Value = $newGenDate

$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""

$wsDe.Range("K2").Value = $epoch

$wsDe.Columns.Item(9).ColumnWidth = 17.8172720046271
$wsDe.Columns.Item(10).ColumnWidth = 20.8721437000093
